$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 9).Value = 'sv'
$ws.Cells.Item(4, 10).Value = 'Statement-opinion'
$ws.Cells.Item(8, 9).Value = 'sv'
$ws.Cells.Item(8, 10).Value = 'Statement-opinion'
$ws.Cells.Item(13, 9).Value = 'ba'
$ws.Cells.Item(13, 10).Value = 'Appreciation'
$ws.Cells.Item(15, 9).Value = 'sd'
$ws.Cells.Item(15, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(17, 9).Value = 'sd'
$ws.Cells.Item(17, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(44, 9).Value = 'sd'
$ws.Cells.Item(44, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(48, 9).Value = 'aa'
$ws.Cells.Item(48, 10).Value = 'Agree/Accept'
$ws.Cells.Item(50, 9).Value = 'b'
$ws.Cells.Item(50, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(68, 9).Value = 'aa'
$ws.Cells.Item(68, 10).Value = 'Agree/Accept'
$ws.Cells.Item(69, 9).Value = 'sd'
$ws.Cells.Item(69, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(70, 9).Value = 'ba'
$ws.Cells.Item(70, 10).Value = 'Appreciation'
$ws.Cells.Item(80, 9).Value = 'aa'
$ws.Cells.Item(80, 10).Value = 'Agree/Accept'
$ws.Cells.Item(98, 9).Value = 'sv'
$ws.Cells.Item(98, 10).Value = 'Statement-opinion'
$ws.Cells.Item(104, 9).Value = 'b'
$ws.Cells.Item(104, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(108, 9).Value = 'b'
$ws.Cells.Item(108, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(115, 9).Value = 'sd'
$ws.Cells.Item(115, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(117, 9).Value = 'sd'
$ws.Cells.Item(117, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(118, 9).Value = 'sv'
$ws.Cells.Item(118, 10).Value = 'Statement-opinion'
$ws.Cells.Item(124, 9).Value = 'b'
$ws.Cells.Item(124, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(125, 9).Value = '%'
$ws.Cells.Item(125, 10).Value = 'Uninterpretable'
$ws.Cells.Item(126, 9).Value = 'sv'
$ws.Cells.Item(126, 10).Value = 'Statement-opinion'
$ws.Cells.Item(128, 9).Value = 'sv'
$ws.Cells.Item(128, 10).Value = 'Statement-opinion'
$ws.Cells.Item(131, 9).Value = 'sv'
$ws.Cells.Item(131, 10).Value = 'Statement-opinion'
$ws.Cells.Item(136, 9).Value = 'sd'
$ws.Cells.Item(136, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(138, 9).Value = 'aa'
$ws.Cells.Item(138, 10).Value = 'Agree/Accept'
$ws.Cells.Item(143, 9).Value = 'sd'
$ws.Cells.Item(143, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(145, 9).Value = 'sd'
$ws.Cells.Item(145, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(148, 9).Value = 'sd'
$ws.Cells.Item(148, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(158, 9).Value = 'b'
$ws.Cells.Item(158, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(159, 9).Value = 'sd'
$ws.Cells.Item(159, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(161, 9).Value = 'ba'
$ws.Cells.Item(161, 10).Value = 'Appreciation'
$ws.Cells.Item(166, 9).Value = 'b'
$ws.Cells.Item(166, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(192, 9).Value = 'sd'
$ws.Cells.Item(192, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(204, 9).Value = 'sd'
$ws.Cells.Item(204, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(205, 9).Value = 'aa'
$ws.Cells.Item(205, 10).Value = 'Agree/Accept'
$ws.Cells.Item(206, 9).Value = 'b'
$ws.Cells.Item(206, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(209, 9).Value = 'sd'
$ws.Cells.Item(209, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(212, 9).Value = 'sd'
$ws.Cells.Item(212, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(225, 9).Value = 'sv'
$ws.Cells.Item(225, 10).Value = 'Statement-opinion'
$ws.Cells.Item(270, 9).Value = 'b'
$ws.Cells.Item(270, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(281, 9).Value = 'sd'
$ws.Cells.Item(281, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(286, 9).Value = 'b'
$ws.Cells.Item(286, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(298, 9).Value = 'sd'
$ws.Cells.Item(298, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(306, 9).Value = 'sv'
$ws.Cells.Item(306, 10).Value = 'Statement-opinion'
$ws.Cells.Item(316, 9).Value = 'b'
$ws.Cells.Item(316, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(322, 9).Value = 'sd'
$ws.Cells.Item(322, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(350, 9).Value = 'aa'
$ws.Cells.Item(350, 10).Value = 'Agree/Accept'
$ws.Cells.Item(360, 9).Value = 'sv'
$ws.Cells.Item(360, 10).Value = 'Statement-opinion'
$ws.Cells.Item(361, 9).Value = 'aa'
$ws.Cells.Item(361, 10).Value = 'Agree/Accept'
$ws.Cells.Item(375, 9).Value = 'aa'
$ws.Cells.Item(375, 10).Value = 'Agree/Accept'
$ws.Cells.Item(376, 9).Value = 'sv'
$ws.Cells.Item(376, 10).Value = 'Statement-opinion'
$ws.Cells.Item(385, 9).Value = 'sv'
$ws.Cells.Item(385, 10).Value = 'Statement-opinion'
$ws.Cells.Item(392, 9).Value = 'b'
$ws.Cells.Item(392, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(400, 9).Value = 'ba'
$ws.Cells.Item(400, 10).Value = 'Appreciation'
$ws.Cells.Item(405, 9).Value = 'sd'
$ws.Cells.Item(405, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(408, 9).Value = 'ba'
$ws.Cells.Item(408, 10).Value = 'Appreciation'
$ws.Cells.Item(421, 9).Value = 'sv'
$ws.Cells.Item(421, 10).Value = 'Statement-opinion'
$ws.Cells.Item(439, 9).Value = 'sv'
$ws.Cells.Item(439, 10).Value = 'Statement-opinion'
$ws.Cells.Item(441, 9).Value = 'ba'
$ws.Cells.Item(441, 10).Value = 'Appreciation'
$ws.Cells.Item(442, 9).Value = 'sd'
$ws.Cells.Item(442, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(443, 9).Value = 'sd'
$ws.Cells.Item(443, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(449, 9).Value = 'b'
$ws.Cells.Item(449, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(452, 9).Value = 'sv'
$ws.Cells.Item(452, 10).Value = 'Statement-opinion'
$ws.Cells.Item(462, 9).Value = 'ba'
$ws.Cells.Item(462, 10).Value = 'Appreciation'
$ws.Cells.Item(468, 9).Value = 'ba'
$ws.Cells.Item(468, 10).Value = 'Appreciation'
$ws.Cells.Item(473, 9).Value = 'b'
$ws.Cells.Item(473, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(478, 9).Value = 'sv'
$ws.Cells.Item(478, 10).Value = 'Statement-opinion'
$ws.Cells.Item(483, 9).Value = 'aa'
$ws.Cells.Item(483, 10).Value = 'Agree/Accept'
$ws.Cells.Item(489, 9).Value = 'ba'
$ws.Cells.Item(489, 10).Value = 'Appreciation'
$ws.Cells.Item(496, 9).Value = 'sd'
$ws.Cells.Item(496, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(503, 9).Value = '%'
$ws.Cells.Item(503, 10).Value = 'Uninterpretable'
$ws.Cells.Item(504, 9).Value = 'sd'
$ws.Cells.Item(504, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(505, 9).Value = 'b'
$ws.Cells.Item(505, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(507, 9).Value = 'sv'
$ws.Cells.Item(507, 10).Value = 'Statement-opinion'
$ws.Cells.Item(523, 9).Value = 'aa'
$ws.Cells.Item(523, 10).Value = 'Agree/Accept'
